$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "iAU_TC_ID_104"
$ws.Range("B2").Value = "@RegressionA Validation of Question Duplicate"
$ws.Range("C2").Value = "failed"
